$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update QUANTITY_IN_STOCK (column D) for rows 2-4 to reflect cart changes
$ws.Range("D2").Value = 2
$ws.Range("D3").Value = 5
$ws.Range("D4").Value = 19
